$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that may look numeric need to be forced to Text format so Excel
# keeps them as strings (matching the source data which stores these as text),
# instead of silently converting them to floating point numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "43.900.84"
$ws.Range("E2").Value = "  +3.00%  "
Set-TextValue $ws.Range("D3") "2.273.31"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  -0.23%  "
Set-TextValue $ws.Range("D5") "321.19"
$ws.Range("E5").Value = "  +2.00%  "
Set-TextValue $ws.Range("D6") "102.99"
$ws.Range("E6").Value = "  +4.49%  "
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +3.05%  "
Set-TextValue $ws.Range("D10") "38.11"
$ws.Range("E10").Value = "  +4.52%  "
Set-TextValue $ws.Range("D11") "0.0844"
$ws.Range("E11").Value = "  +2.75%  "
Set-TextValue $ws.Range("D12") "7.86"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  +3.05%  "
Set-TextValue $ws.Range("D14") "2.614.78"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("E15").Value = "  +3.15%  "
Set-TextValue $ws.Range("D16") "14.57"
$ws.Range("E16").Value = "  +4.16%  "
Set-TextValue $ws.Range("D17") "2.271.20"
$ws.Range("E17").Value = "  +3.43%  "
Set-TextValue $ws.Range("D18") "43.856.20"
$ws.Range("E18").Value = "  +3.23%  "
Set-TextValue $ws.Range("D19") "14.32"
$ws.Range("E19").Value = "  +3.42%  "
Set-TextValue $ws.Range("D20") "0.0₃0995"
$ws.Range("E20").Value = "  +4.26%  "
Set-TextValue $ws.Range("D21") "6.70"
$ws.Range("E21").Value = "  +3.49%  "
Set-TextValue $ws.Range("D22") "66.16"
$ws.Range("E22").Value = "  +0.79%  "
Set-TextValue $ws.Range("D23") "3.22"
$ws.Range("E23").Value = "  +1.28%  "
Set-TextValue $ws.Range("D24") "239.36"
$ws.Range("E24").Value = "  +2.22%  "
Set-TextValue $ws.Range("D25") "2.23"
$ws.Range("E25").Value = "  +4.88%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +3.14%  "
Set-TextValue $ws.Range("D28") "10.21"
$ws.Range("E28").Value = "  +1.31%  "
Set-TextValue $ws.Range("D29") "39.26"
$ws.Range("E29").Value = "  +16.63%  "
Set-TextValue $ws.Range("D30") "2.20"
$ws.Range("E30").Value = "  +1.89%  "
Set-TextValue $ws.Range("D31") "6.50"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("E33").Value = "  +0.57%  "
Set-TextValue $ws.Range("D34") "161.18"
$ws.Range("E34").Value = "  +2.16%  "
Set-TextValue $ws.Range("D35") "3.43"
$ws.Range("E35").Value = "  +7.28%  "
$ws.Range("E36").Value = "  +0.82%  "
Set-TextValue $ws.Range("D37") "2.04"
$ws.Range("E37").Value = "  +10.11%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D40") "3.92"
$ws.Range("E40").Value = "  +12.12%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.107"
$ws.Range("E41").Value = "  +4.29%  "
Set-TextValue $ws.Range("D42") "15.72"
$ws.Range("E42").Value = "  +32.17%  "
$ws.Range("E43").Value = "  +3.51%  "
$ws.Range("E44").Value = "  -0.11%  "
Set-TextValue $ws.Range("D45") "1.819.64"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("E46").Value = "  +1.96%  "
Set-TextValue $ws.Range("D47") "86.30"
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("E48").Value = "  +2.16%  "
Set-TextValue $ws.Range("D49") "76.71"
$ws.Range("E49").Value = "  -0.08%  "
Set-TextValue $ws.Range("D50") "8.87"
$ws.Range("E50").Value = "  +5.53%  "
Set-TextValue $ws.Range("D51") "59.82"
$ws.Range("E51").Value = "  -0.08%  "
